$d = $word.ActiveDocument

$replacements = @(
    @("933×3=2799", "361×5=1805"),
    @("234×9=2106", "358×9=3222"),
    @("137×3=411",  "382×6=2292"),
    @("335×7=2345", "732×4=2928"),
    @("479×5=2395", "577×8=4616"),
    @("550×9=4950", "685×3=2055"),
    @("375×6=2250", "219×7=1533"),
    @("680×5=3400", "558×9=5022"),
    @("275×4=1100", "195×8=1560"),
    @("694×7=4858", "443×3=1329"),
    @("883×6=5298", "737×6=4422"),
    @("476×9=4284", "998×8=7984"),
    @("716×4=2864", "164×2=328"),
    @("794×5=3970", "795×7=5565"),
    @("428×7=2996", "536×4=2144"),
    @("843×9=7587", "673×2=1346"),
    @("649×8=5192", "246×9=2214"),
    @("120×6=720",  "723×9=6507"),
    @("514×8=4112", "112×7=784"),
    @("564×2=1128", "251×5=1255"),
    @("424×8=3392", "149×2=298"),
    @("951×5=4755", "520×9=4680"),
    @("423×2=846",  "746×7=5222"),
    @("144×7=1008", "906×7=6342"),
    @("692×2=1384", "118×2=236")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
